$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.834.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.559.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.09"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.49"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.65"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.42"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.950.03"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.98"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.580.84"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.843.16"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0959"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.86"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.39"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.85"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.21"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.76"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0796"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.29"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.65"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.71"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.59"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +12.14%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.79"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.07"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.81%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.988.55"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.95"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.804.15"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.39"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.194"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.48"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.06%  "
